$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.309.96"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "2.345.89"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "2.344.36"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "2.762.98"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "60.252.34"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "2.359.27"
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  +7.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "314.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +8.31%  "
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +12.34%  "
$ws.Range("D32").Value = "0.0₃0729"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.382"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +8.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "324.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.14%  "
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +9.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("D50").Value = "0.0₆0211"
$ws.Range("E50").Value = "  +15.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.99%  "
